# --------------------------------------------------------------------------
# Helper: push a list of "RRGGBB" hex strings into a ThemeColorScheme's
# 12 slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), in order.
# PowerPoint's ColorFormat.RGB / ThemeColor.RGB values are packed as
# 0xBBGGRR (the usual OLE_COLOR convention), so we repack before assigning.
# --------------------------------------------------------------------------
function Set-ThemeColors {
    param($ThemeColorScheme, [string[]]$HexColors)

    for ($i = 0; $i -lt $HexColors.Length; $i++) {
        $hex = $HexColors[$i]
        $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
        $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
        $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
        $bgr = ($b * 65536) + ($g * 256) + $r
        $ThemeColorScheme.Colors($i + 1).RGB = $bgr
    }
}

$p = $ppt.ActivePresentation

# --------------------------------------------------------------------------
# 1) Slide 16's table switches from the deck's one custom table style to a
#    built-in PowerPoint table style.
# --------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{155E508A-D569-486B-9768-0C488B1D81E0}")
    }
}

# --------------------------------------------------------------------------
# 2) The deck carries two theme colour schemes: the main design ("Integral")
#    used by the slide master, and a default "Office Theme" used only by the
#    notes master. Their colour schemes are swapped - the slide master now
#    takes on the Office Theme palette and the notes master takes on the
#    Integral palette.
# --------------------------------------------------------------------------
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)
$integralColors = @(
    "000000", "FFFFFF", "455F51", "E3DED1",
    "99CB38", "63A537", "E6D024", "CC9700",
    "4EB3CF", "378DA6", "6B9F25", "B26B02"
)

$slideMasterThemeColors = $p.SlideMaster.Theme.ThemeColorScheme
Set-ThemeColors $slideMasterThemeColors $officeColors

$notesMasterThemeColors = $p.NotesMaster.Theme.ThemeColorScheme
Set-ThemeColors $notesMasterThemeColors $integralColors
